# Scheduled runner update: refresh market-price / profit columns (H..N)
# on the Pandaemonium_Profits leve-crafting sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) with newly polled Universalis price data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 197.33333
$ws.Range("I6").Value = 138.33333
$ws.Range("K6").Value = 414.99999
$ws.Range("M6").Value = -302.99999

$ws.Range("H98").Value = 1354.0358
$ws.Range("I98").Value = 909.26086
$ws.Range("J98").Value = 3400
$ws.Range("K98").Value = 909.26086
$ws.Range("L98").Value = 3400
$ws.Range("M98").Value = 588.73914
$ws.Range("N98").Value = -6396

$ws.Range("H103").Value = 1255.6154
$ws.Range("I103").Value = 1300
$ws.Range("J103").Value = 1251.9166
$ws.Range("K103").Value = 3900
$ws.Range("L103").Value = 3755.7498
$ws.Range("M103").Value = -3314
$ws.Range("N103").Value = -4927.7498

$ws.Range("H122").Value = 1354.0358
$ws.Range("I122").Value = 909.26086
$ws.Range("J122").Value = 3400
$ws.Range("K122").Value = 2727.78258
$ws.Range("L122").Value = 10200
$ws.Range("M122").Value = -277.7825800000001
$ws.Range("N122").Value = -15100

$ws.Range("H126").Value = 57445
$ws.Range("J126").Value = 57445
$ws.Range("L126").Value = 57445
$ws.Range("N126").Value = -67325

$ws.Range("H129").Value = 955.4407
$ws.Range("I129").Value = 324.4
$ws.Range("J129").Value = 1013.87036
$ws.Range("K129").Value = 973.1999999999999
$ws.Range("L129").Value = 3041.61108
$ws.Range("M129").Value = 4026.8
$ws.Range("N129").Value = -13041.61108

$ws.Range("H137").Value = 3161.2942
$ws.Range("I137").Value = 3095.8572
$ws.Range("J137").Value = 3466.6667
$ws.Range("K137").Value = 9287.571599999999
$ws.Range("L137").Value = 10400.0001
$ws.Range("M137").Value = -6737.571599999999
$ws.Range("N137").Value = -15500.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2275503
$ws.Range("I2").Value = 2275503
$ws.Range("K2").Value = 2275503
$ws.Range("M2").Value = -2275390

$ws.Range("H74").Value = 4199.756
$ws.Range("I74").Value = 1804.7858
$ws.Range("J74").Value = 9358.154
$ws.Range("K74").Value = 1804.7858
$ws.Range("L74").Value = 9358.154
$ws.Range("M74").Value = -930.7858000000001
$ws.Range("N74").Value = -11106.154

$ws.Range("H77").Value = 4199.756
$ws.Range("I77").Value = 1804.7858
$ws.Range("J77").Value = 9358.154
$ws.Range("K77").Value = 9023.929
$ws.Range("L77").Value = 46790.77
$ws.Range("M77").Value = -4655.929
$ws.Range("N77").Value = -55526.77

$ws.Range("H110").Value = 1219.5
$ws.Range("I110").Value = 1174.1333
$ws.Range("K110").Value = 1174.1333
$ws.Range("M110").Value = 870.8667

$ws.Range("H116").Value = 2275503
$ws.Range("I116").Value = 2275503
$ws.Range("K116").Value = 2275503
$ws.Range("M116").Value = -2273209

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2275503
$ws.Range("I3").Value = 2275503
$ws.Range("K3").Value = 2275503
$ws.Range("M3").Value = -2275389

$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -127
$ws.Range("N22").ClearContents()

$ws.Range("H105").Value = 2088094.1
$ws.Range("I105").Value = 2844591
$ws.Range("J105").Value = 7727.75
$ws.Range("K105").Value = 2844591
$ws.Range("L105").Value = 7727.75
$ws.Range("M105").Value = -2842844
$ws.Range("N105").Value = -11221.75

$ws.Range("H107").Value = 1225.9062
$ws.Range("I107").Value = 1158.15
$ws.Range("J107").Value = 1338.8334
$ws.Range("K107").Value = 1158.15
$ws.Range("L107").Value = 1338.8334
$ws.Range("M107").Value = 761.8499999999999
$ws.Range("N107").Value = -5178.8334

$ws.Range("H141").Value = 41693.332
$ws.Range("J141").Value = 41693.332
$ws.Range("L141").Value = 41693.332
$ws.Range("N141").Value = -52053.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4023.8823
$ws.Range("I31").Value = 3965.5557
$ws.Range("J31").Value = 4163.8667
$ws.Range("K31").Value = 3965.5557
$ws.Range("L31").Value = 4163.8667
$ws.Range("M31").Value = -3670.5557
$ws.Range("N31").Value = -4753.8667

$ws.Range("H34").Value = 4023.8823
$ws.Range("I34").Value = 3965.5557
$ws.Range("J34").Value = 4163.8667
$ws.Range("K34").Value = 3965.5557
$ws.Range("L34").Value = 4163.8667
$ws.Range("M34").Value = -3763.5557
$ws.Range("N34").Value = -4567.8667

$ws.Range("H86").Value = 2490.6365
$ws.Range("I86").Value = 2589.6667
$ws.Range("J86").Value = 2045
$ws.Range("K86").Value = 2589.6667
$ws.Range("L86").Value = 2045
$ws.Range("M86").Value = -1466.6667
$ws.Range("N86").Value = -4291

$ws.Range("H89").Value = 2490.6365
$ws.Range("I89").Value = 2589.6667
$ws.Range("J89").Value = 2045
$ws.Range("K89").Value = 12948.3335
$ws.Range("L89").Value = 10225
$ws.Range("M89").Value = -7332.333500000001
$ws.Range("N89").Value = -21457

$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 29.153847
$ws.Range("I2").Value = 22.285715
$ws.Range("J2").Value = 30.65625
$ws.Range("K2").Value = 133.71429
$ws.Range("L2").Value = 183.9375
$ws.Range("M2").Value = -20.71429000000001
$ws.Range("N2").Value = -409.9375

$ws.Range("H46").Value = 2933.4219
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 2964.111
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 8892.332999999999
$ws.Range("M46").Value = -2909
$ws.Range("N46").Value = -9074.332999999999

$ws.Range("H114").Value = 1133.3334
$ws.Range("I114").Value = 700
$ws.Range("K114").Value = 2100
$ws.Range("M114").Value = 1154

$ws.Range("H123").Value = 2000
$ws.Range("J123").Value = 3000
$ws.Range("L123").Value = 9000
$ws.Range("N123").Value = -13900

$ws.Range("H137").Value = 32720.941
$ws.Range("I137").Value = 1148.3334
$ws.Range("J137").Value = 49942.363
$ws.Range("K137").Value = 3445.0002
$ws.Range("L137").Value = 149827.089
$ws.Range("M137").Value = 1654.9998
$ws.Range("N137").Value = -160027.089

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6574.636
$ws.Range("I122").Value = 10102.167
$ws.Range("J122").Value = 2341.6
$ws.Range("K122").Value = 30306.501
$ws.Range("L122").Value = 7024.799999999999
$ws.Range("M122").Value = -27856.501
$ws.Range("N122").Value = -11924.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6170.913
$ws.Range("I122").Value = 5615.7188
$ws.Range("K122").Value = 16847.1564
$ws.Range("M122").Value = -14397.1564

$ws.Range("H132").Value = 3047.9546
$ws.Range("I132").Value = 2629.1853
$ws.Range("J132").Value = 3713.0588
$ws.Range("K132").Value = 7887.5559
$ws.Range("L132").Value = 11139.1764
$ws.Range("M132").Value = -5357.5559
$ws.Range("N132").Value = -16199.1764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 12059.7
$ws.Range("I122").Value = 1941.7142
$ws.Range("J122").Value = 35668.332
$ws.Range("K122").Value = 5825.142599999999
$ws.Range("L122").Value = 107004.996
$ws.Range("M122").Value = -3375.142599999999
$ws.Range("N122").Value = -111904.996

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
